$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,2).Value = 1.02
$ws.Cells.Item(2,3).Value = 1.027948233614345
$ws.Cells.Item(2,4).Value = 1.033086633838193
$ws.Cells.Item(2,5).Value = 1.031580492013402
$ws.Cells.Item(2,6).Value = 1.03846161479219
$ws.Cells.Item(2,9).Value = 1.035017703491167
$ws.Cells.Item(2,10).Value = 1.033103271979283
$ws.Cells.Item(2,11).Value = 1.03588999523644
$ws.Cells.Item(2,12).Value = 1.034388198615864
$ws.Cells.Item(2,13).Value = 1.041249581012584
$ws.Cells.Item(2,14).Value = 1.015040374564616

$ws.Cells.Item(3,2).Value = 1.02
$ws.Cells.Item(3,3).Value = 1.028790566802157
$ws.Cells.Item(3,4).Value = 1.033723155003162
$ws.Cells.Item(3,5).Value = 1.032370208662422
$ws.Cells.Item(3,6).Value = 1.039562938878211
$ws.Cells.Item(3,9).Value = 1.035219498487311
$ws.Cells.Item(3,10).Value = 1.033586556454037
$ws.Cells.Item(3,11).Value = 1.036335884269723
$ws.Cells.Item(3,12).Value = 1.034986559326036
$ws.Cells.Item(3,13).Value = 1.042160153827706
$ws.Cells.Item(3,14).Value = 1.015202268941144

$ws.Cells.Item(4,2).Value = 1.02
$ws.Cells.Item(4,3).Value = 1.029336116413207
$ws.Cells.Item(4,4).Value = 1.034135418489456
$ws.Cells.Item(4,5).Value = 1.032882055001304
$ws.Cells.Item(4,6).Value = 1.04027645184345
$ws.Cells.Item(4,9).Value = 1.035349142042703
$ws.Cells.Item(4,10).Value = 1.033899129518568
$ws.Cells.Item(4,11).Value = 1.036624107686115
$ws.Cells.Item(4,12).Value = 1.035373930319156
$ws.Cells.Item(4,13).Value = 1.042749649721957
$ws.Cells.Item(4,14).Value = 1.015306929609477

$ws.Cells.Item(5,2).Value = 1.02
$ws.Cells.Item(5,3).Value = 1.029565584741027
$ws.Cells.Item(5,4).Value = 1.034308826327321
$ws.Cells.Item(5,5).Value = 1.033097436376912
$ws.Cells.Item(5,6).Value = 1.040576623280653
$ws.Cells.Item(5,9).Value = 1.035403420668078
$ws.Cells.Item(5,10).Value = 1.034030499682421
$ws.Cells.Item(5,11).Value = 1.03674520476178
$ws.Cells.Item(5,12).Value = 1.035536825727539
$ws.Cells.Item(5,13).Value = 1.042997543531492
$ws.Cells.Item(5,14).Value = 1.015350905626157

$ws.Cells.Item(6,2).Value = 1.02
$ws.Cells.Item(6,3).Value = 1.029604120438842
$ws.Cells.Item(6,4).Value = 1.034337947632975
$ws.Cells.Item(6,5).Value = 1.033133611609078
$ws.Cells.Item(6,6).Value = 1.040627035713949
$ws.Cells.Item(6,9).Value = 1.035412521163955
$ws.Cells.Item(6,10).Value = 1.034052555196971
$ws.Cells.Item(6,11).Value = 1.036765533238419
$ws.Cells.Item(6,12).Value = 1.035564179165596
$ws.Cells.Item(6,13).Value = 1.043039170060762
$ws.Cells.Item(6,14).Value = 1.015358288013566

$ws.Cells.Item(7,2).Value = 1.02
$ws.Cells.Item(7,3).Value = 1.029339182113223
$ws.Cells.Item(7,4).Value = 1.034137735212822
$ws.Cells.Item(7,5).Value = 1.03288493214927
$ws.Cells.Item(7,6).Value = 1.040280461922307
$ws.Cells.Item(7,9).Value = 1.03534986819536
$ws.Cells.Item(7,10).Value = 1.033900885034081
$ws.Cells.Item(7,11).Value = 1.036625726075572
$ws.Cells.Item(7,12).Value = 1.035376106761696
$ws.Cells.Item(7,13).Value = 1.042752961818078
$ws.Cells.Item(7,14).Value = 1.015307517311316

$ws.Cells.Item(8,2).Value = 1.019999999999999
$ws.Cells.Item(8,3).Value = 1.028232798967984
$ws.Cells.Item(8,4).Value = 1.033301667126828
$ws.Cells.Item(8,5).Value = 1.031847204163252
$ws.Cells.Item(8,6).Value = 1.038833629233082
$ws.Cells.Item(8,9).Value = 1.035086093543548
$ws.Cells.Item(8,10).Value = 1.033266629505392
$ws.Cells.Item(8,11).Value = 1.036040746310613
$ws.Cells.Item(8,12).Value = 1.034590377008866
$ws.Cells.Item(8,13).Value = 1.041557251953771
$ws.Cells.Item(8,14).Value = 1.015095107133197

$ws.Cells.Item(9,2).Value = 1.02
$ws.Cells.Item(9,3).Value = 1.02628712667078
$ws.Cells.Item(9,4).Value = 1.031831475309607
$ws.Cells.Item(9,5).Value = 1.030025149963794
$ws.Cells.Item(9,6).Value = 1.036290916891819
$ws.Cells.Item(9,9).Value = 1.034614181487251
$ws.Cells.Item(9,10).Value = 1.03214792941562
$ws.Cells.Item(9,11).Value = 1.035007711583941
$ws.Cells.Item(9,12).Value = 1.033207342457132
$ws.Cells.Item(9,13).Value = 1.039452545953176
$ws.Cells.Item(9,14).Value = 1.014720096578213

$ws.Cells.Item(10,2).Value = 1.02
$ws.Cells.Item(10,3).Value = 1.024992722809713
$ws.Cells.Item(10,4).Value = 1.030853499355124
$ws.Cells.Item(10,5).Value = 1.028814947390779
$ws.Cells.Item(10,6).Value = 1.034600382814753
$ws.Cells.Item(10,9).Value = 1.034294828648616
$ws.Cells.Item(10,10).Value = 1.031401475660446
$ws.Cells.Item(10,11).Value = 1.034317587110315
$ws.Cells.Item(10,12).Value = 1.032286415275446
$ws.Cells.Item(10,13).Value = 1.038050989520949
$ws.Cells.Item(10,14).Value = 1.014469630736848

$ws.Cells.Item(11,2).Value = 1.02
$ws.Cells.Item(11,3).Value = 1.02443289090037
$ws.Cells.Item(11,4).Value = 1.030430554275738
$ws.Cells.Item(11,5).Value = 1.02829200322271
$ws.Cells.Item(11,6).Value = 1.033869463836392
$ws.Cells.Item(11,9).Value = 1.034155426993631
$ws.Cells.Item(11,10).Value = 1.031078110336627
$ws.Cells.Item(11,11).Value = 1.034018428239962
$ws.Cells.Item(11,12).Value = 1.031887918801428
$ws.Cells.Item(11,13).Value = 1.037444483655509
$ws.Cells.Item(11,14).Value = 1.014361072266182

$ws.Cells.Item(12,2).Value = 1.02
$ws.Cells.Item(12,3).Value = 1.024225043866185
$ws.Cells.Item(12,4).Value = 1.030273534030371
$ws.Cells.Item(12,5).Value = 1.028097922529816
$ws.Cells.Item(12,6).Value = 1.033598132844049
$ws.Cells.Item(12,9).Value = 1.034103479324918
$ws.Cells.Item(12,10).Value = 1.030957977118473
$ws.Cells.Item(12,11).Value = 1.033907258698284
$ws.Cells.Item(12,12).Value = 1.03173994140773
$ws.Cells.Item(12,13).Value = 1.037219257688955
$ws.Cells.Item(12,14).Value = 1.014320733408224

$ws.Cells.Item(13,2).Value = 1.02
$ws.Cells.Item(13,3).Value = 1.024269623239025
$ws.Cells.Item(13,4).Value = 1.030307211743892
$ws.Cells.Item(13,5).Value = 1.028139546048262
$ws.Cells.Item(13,6).Value = 1.033656326793219
$ws.Cells.Item(13,9).Value = 1.034114629855619
$ws.Cells.Item(13,10).Value = 1.030983747033298
$ws.Cells.Item(13,11).Value = 1.033931107125925
$ws.Cells.Item(13,12).Value = 1.031771681143637
$ws.Cells.Item(13,13).Value = 1.037267566819471
$ws.Cells.Item(13,14).Value = 1.014329386921305

$ws.Cells.Item(14,2).Value = 1.02
$ws.Cells.Item(14,3).Value = 1.024415708163235
$ws.Cells.Item(14,4).Value = 1.030417573284589
$ws.Cells.Item(14,5).Value = 1.028275957096353
$ws.Cells.Item(14,6).Value = 1.033847032165412
$ws.Cells.Item(14,9).Value = 1.034151136400735
$ws.Cells.Item(14,10).Value = 1.031068180508795
$ws.Cells.Item(14,11).Value = 1.034009239916817
$ws.Cells.Item(14,12).Value = 1.03187568608132
$ws.Cells.Item(14,13).Value = 1.037425865233953
$ws.Cells.Item(14,14).Value = 1.014357738155712

$ws.Cells.Item(15,2).Value = 1.02
$ws.Cells.Item(15,3).Value = 1.024505729124731
$ws.Cells.Item(15,4).Value = 1.030485581369895
$ws.Cells.Item(15,5).Value = 1.028360026230994
$ws.Cells.Item(15,6).Value = 1.033964553902514
$ws.Cells.Item(15,9).Value = 1.034173607084221
$ws.Cells.Item(15,10).Value = 1.031120200014528
$ws.Cells.Item(15,11).Value = 1.034057373696398
$ws.Cells.Item(15,12).Value = 1.031939772549779
$ws.Cells.Item(15,13).Value = 1.037523405725546
$ws.Cells.Item(15,14).Value = 1.01437520425663

$ws.Cells.Item(16,2).Value = 1.02
$ws.Cells.Item(16,3).Value = 1.025029890887694
$ws.Cells.Item(16,4).Value = 1.030881580017176
$ws.Cells.Item(16,5).Value = 1.028849676393366
$ws.Cells.Item(16,6).Value = 1.034648914587208
$ws.Cells.Item(16,9).Value = 1.034304056731616
$ws.Cells.Item(16,10).Value = 1.031422933361584
$ws.Cells.Item(16,11).Value = 1.03433743442518
$ws.Cells.Item(16,12).Value = 1.032312867988837
$ws.Cells.Item(16,13).Value = 1.038091249320939
$ws.Cells.Item(16,14).Value = 1.014476833223197

$ws.Cells.Item(17,2).Value = 1.02
$ws.Cells.Item(17,3).Value = 1.025358859578906
$ws.Cells.Item(17,4).Value = 1.031130121126863
$ws.Cells.Item(17,5).Value = 1.02915711166722
$ws.Cells.Item(17,6).Value = 1.035078489117463
$ws.Cells.Item(17,9).Value = 1.034385584939504
$ws.Cells.Item(17,10).Value = 1.031612791619099
$ws.Cells.Item(17,11).Value = 1.034513021381833
$ws.Cells.Item(17,12).Value = 1.032546974331655
$ws.Cells.Item(17,13).Value = 1.038447544066824
$ws.Cells.Item(17,14).Value = 1.014540554518823

$ws.Cells.Item(18,2).Value = 1.02
$ws.Cells.Item(18,3).Value = 1.025550804418328
$ws.Cells.Item(18,4).Value = 1.03127514150068
$ws.Cells.Item(18,5).Value = 1.029336537708167
$ws.Cells.Item(18,6).Value = 1.035329158133938
$ws.Cells.Item(18,9).Value = 1.034433030885711
$ws.Cells.Item(18,10).Value = 1.031723518583175
$ws.Cells.Item(18,11).Value = 1.034615406333573
$ws.Cells.Item(18,12).Value = 1.032683550673818
$ws.Cells.Item(18,13).Value = 1.038655401127317
$ws.Cells.Item(18,14).Value = 1.014577711905365

$ws.Cells.Item(19,2).Value = 1.02
$ws.Cells.Item(19,3).Value = 1.025616263289709
$ws.Cells.Item(19,4).Value = 1.031324598228318
$ws.Cells.Item(19,5).Value = 1.029397734985773
$ws.Cells.Item(19,6).Value = 1.035414647661001
$ws.Cells.Item(19,9).Value = 1.034449190378772
$ws.Cells.Item(19,10).Value = 1.031761271187251
$ws.Cells.Item(19,11).Value = 1.034650311505916
$ws.Cells.Item(19,12).Value = 1.032730124062571
$ws.Cells.Item(19,13).Value = 1.038726281162246
$ws.Cells.Item(19,14).Value = 1.014590379877768

$ws.Cells.Item(20,2).Value = 1.02
$ws.Cells.Item(20,3).Value = 1.025323557829372
$ws.Cells.Item(20,4).Value = 1.031103449778418
$ws.Cells.Item(20,5).Value = 1.029124115977788
$ws.Cells.Item(20,6).Value = 1.035032388923236
$ws.Cells.Item(20,9).Value = 1.03437684890946
$ws.Cells.Item(20,10).Value = 1.03159242308088
$ws.Cells.Item(20,11).Value = 1.034494185860394
$ws.Cells.Item(20,12).Value = 1.032521854230989
$ws.Cells.Item(20,13).Value = 1.038409313246723
$ws.Cells.Item(20,14).Value = 1.01453371887572

$ws.Cells.Item(21,2).Value = 1.02
$ws.Cells.Item(21,3).Value = 1.024372687038137
$ws.Cells.Item(21,4).Value = 1.030385072330231
$ws.Cells.Item(21,5).Value = 1.028235782901334
$ws.Cells.Item(21,6).Value = 1.033790869630822
$ws.Cells.Item(21,9).Value = 1.034140390757698
$ws.Cells.Item(21,10).Value = 1.031043317524823
$ws.Cells.Item(21,11).Value = 1.03398623309115
$ws.Cells.Item(21,12).Value = 1.031845058051931
$ws.Cells.Item(21,13).Value = 1.037379248709073
$ws.Cells.Item(21,14).Value = 1.014349389846087

$ws.Cells.Item(22,2).Value = 1.02
$ws.Cells.Item(22,3).Value = 1.023775412933757
$ws.Cells.Item(22,4).Value = 1.029933865771209
$ws.Cells.Item(22,5).Value = 1.02767820258673
$ws.Cells.Item(22,6).Value = 1.033011231302305
$ws.Cells.Item(22,9).Value = 1.03399075021784
$ws.Cells.Item(22,10).Value = 1.030697952408642
$ws.Cells.Item(22,11).Value = 1.03366658224796
$ws.Cells.Item(22,12).Value = 1.031419772684268
$ws.Cells.Item(22,13).Value = 1.036731938183092
$ws.Cells.Item(22,14).Value = 1.014233405910404

$ws.Cells.Item(23,2).Value = 1.02
$ws.Cells.Item(23,3).Value = 1.024091984200524
$ws.Cells.Item(23,4).Value = 1.03017301424089
$ws.Cells.Item(23,5).Value = 1.027973695893923
$ws.Cells.Item(23,6).Value = 1.033424441694167
$ws.Cells.Item(23,9).Value = 1.034070169253362
$ws.Cells.Item(23,10).Value = 1.030881048099544
$ws.Cells.Item(23,11).Value = 1.033836061480099
$ws.Cells.Item(23,12).Value = 1.03164520104133
$ws.Cells.Item(23,13).Value = 1.03707505806545
$ws.Cells.Item(23,14).Value = 1.014294899493638

$ws.Cells.Item(24,2).Value = 1.02
$ws.Cells.Item(24,3).Value = 1.025339508976013
$ws.Cells.Item(24,4).Value = 1.031115501256719
$ws.Cells.Item(24,5).Value = 1.029139024987639
$ws.Cells.Item(24,6).Value = 1.035053219290771
$ws.Cells.Item(24,9).Value = 1.034380796679398
$ws.Cells.Item(24,10).Value = 1.031601626790361
$ws.Cells.Item(24,11).Value = 1.034502696920599
$ws.Cells.Item(24,12).Value = 1.032533204842773
$ws.Cells.Item(24,13).Value = 1.038426587996342
$ws.Cells.Item(24,14).Value = 1.01453680764006

$ws.Cells.Item(25,2).Value = 1.02
$ws.Cells.Item(25,3).Value = 1.026789657486778
$ws.Cells.Item(25,4).Value = 1.032211182416001
$ws.Cells.Item(25,5).Value = 1.030495408257131
$ws.Cells.Item(25,6).Value = 1.036947459916986
$ws.Cells.Item(25,9).Value = 1.034737020853355
$ws.Cells.Item(25,10).Value = 1.032437259783655
$ws.Cells.Item(25,11).Value = 1.035275033211854
$ws.Cells.Item(25,12).Value = 1.033564702217594
$ws.Cells.Item(25,13).Value = 1.039996387421399
$ws.Cells.Item(25,14).Value = 1.01481712835003
